# bulk_user_details.xlsx — add new iAuthor testcases
# Row 2 gets a refreshed set of test credentials, and two brand-new
# candidate rows (3 and 4) are appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- give the two new rows the same thin-border cell style used by row 2 ---
$ws.Range("A3:H4").Borders.LineStyle = 1

# --- row 2: replace the existing test data values ---
$ws.Range("A2").Value = "eilHw781"
$ws.Range("C2").Value = "stubfxt98"
$ws.Range("D2").Value = "PkT$26!x"
$ws.Range("F2").Value = "HqGPEIfK"
$ws.Range("G2").Value = "Zkpf"
$ws.Range("B2").Value = 23101743

# --- row 3: new candidate ---
$ws.Range("A3").Value = "kukFV955"
$ws.Range("B3").Value = 23101742
$ws.Range("C3").Value = "wfsbkeg56"
$ws.Range("D3").Value = "N9&p6k#Z"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "ZudVdWmp"
$ws.Range("G3").Value = "GvQr"
$ws.Range("H3").Value = "Candidate"

# --- row 4: new candidate ---
$ws.Range("A4").Value = "HLknW522"
$ws.Range("B4").Value = 23101741
$ws.Range("C4").Value = "mjjhtkq74"
$ws.Range("D4").Value = "GSq4&3!w"
$ws.Range("E4").Value = "MR"
$ws.Range("F4").Value = "liQjBovl"
$ws.Range("G4").Value = "yNkN"
$ws.Range("H4").Value = "Candidate"

# --- the used range now spans through row 4; refresh the selection to match ---
$ws.Range("A1:H4").Select() | Out-Null
